$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 26.05.2025"

$ws.Range("B6").Value = "27.05."
$ws.Range("C6").Value = "28.05."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU GIWFHN"
$ws.Range("E6").Value = "162,57-"

$ws.Range("B7").Value = "30.05."
$ws.Range("C7").Value = "31.05."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU BBUPRH"
$ws.Range("E7").Value = "180,76-"

$ws.Range("B8").Value = "31.05."
$ws.Range("C8").Value = "01.06."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 40097606"
$ws.Range("E8").Value = "86,55-"

$ws.Range("D12").Value = "KONTOSTAND AM 05.06.2025"
$ws.Range("E12").Value = "429,88-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.06.2025"
